$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: J2 becomes the "лаб1" label (was a blank placeholder string),
#    K2 ("тк" header) loses its text but keeps its style.
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "лаб1"
$ws.Range("K2").ClearContents()

# ---------------------------------------------------------------------------
# 2. Drop the old per-row "total / count / flag" helper columns J:L.
#    J4:J31 held =SUM(...) formulas, K and L held plain numbers - all of it
#    goes away except for the two rows (6 and 26) that get a literal lab1
#    score typed directly into J, handled below.
# ---------------------------------------------------------------------------
$ws.Range("J4:J31").ClearContents()
$ws.Range("K1:L32").ClearContents()

# ---------------------------------------------------------------------------
# 3. Row 6 (student #3): fill in the rest of the grades and the lab1 score.
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 5

# copy the plain thick-border style (already used on G/H of this row) onto F6
$ws.Range("G6").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
# I6/J6 need the thin "edge" style used elsewhere for the last filled column
# (same style as I26, taken from that existing cell)
$ws.Range("I26").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 4. Row 14 (student #11): bump grades to 5 and restyle C14:F14 to the
#    green "inner" look (same style already used on C9/C10 etc. minus edges).
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 5

$ws.Range("I24").Copy() | Out-Null
$ws.Range("C14:F14").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 5. Row 26 (student #23): literal lab1 score instead of the old SUM formula.
# ---------------------------------------------------------------------------
$ws.Range("J26").Value = 5
$ws.Range("I26").Copy() | Out-Null
$ws.Range("J26").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 6. New conditional formatting (3-colour scale) on D14, on top of / ahead of
#    the existing one covering J4:J31.
# ---------------------------------------------------------------------------
$newCf = $ws.Range("D14").FormatConditions.AddColorScale(3)
$newCf.SetFirstPriority() | Out-Null

# ---------------------------------------------------------------------------
# 7. Selection, as left by the editor.
# ---------------------------------------------------------------------------
$ws.Range("J6").Select() | Out-Null
